$d = $word.ActiveDocument

$replacements = @(
    @("33×47=", "95×94="),
    @("13×39=", "83×43="),
    @("42×50=", "11×50="),
    @("87×51=", "83×22="),
    @("59×79=", "97×85="),
    @("23×29=", "52×73="),
    @("34×24=", "75×92="),
    @("13×42=", "12×97="),
    @("55×65=", "11×90="),
    @("92×11=", "29×68="),
    @("47×14=", "97×94="),
    @("83×55=", "52×23="),
    @("19×28=", "29×83="),
    @("16×84=", "21×14="),
    @("95×57=", "65×98="),
    @("95×35=", "13×11="),
    @("14×57=", "39×61="),
    @("93×22=", "50×32="),
    @("14×96=", "53×31="),
    @("39×37=", "68×14="),
    @("34×12=", "22×85="),
    @("23×40=", "97×76="),
    @("67×30=", "87×39="),
    @("28×22=", "82×46="),
    @("65×49=", "66×81=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
